# Legobot turn data update:
#  - fill in the missing "Power 80" row values (K4, J5, K5) to match the
#    pattern already present for the other power rows
#  - move the active selection to F41

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K4").Value = 30
$ws.Range("J5").Value = 42.9
$ws.Range("K5").Value = 26.7

$ws.Range("F41").Select() | Out-Null
